$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set for rows 2..17 (columns A..K)
# A, B, C, D, E, F, G, H, I, J, K
$data = @(
    @(9,  10, 340.6121553362682,  0.02222224,          0.02222224,          0.04444448, 510.9182330044023, 0.38, 0.47, 0.18, 0.39),
    @(10, 11, 428.3635523572899,  0.08611118000000001, 0.08611118000000001, 0.17222236, 642.5453285359349, 0.4,  0.45, 0.12, 0.4),
    @(11, 12, 498.7069374946091,  0.10000008,          0.10000008,          0.20000016, 748.0604062419137, 0.34, 0.42, 0.14, 0.3),
    @(12, 13, 508.3495138842192,  0.11944454,          0.11944454,          0.23888908, 762.5242708263288, 0.35, 0.48, 0.12, 0.36),
    @(13, 14, 425.5577207029365,  0.10277786,          0.10277786,          0.20555572, 638.3365810544047, 0.35, 0.45, 0.14, 0.38),
    @(14, 15, 334.7633794933903,  0.0416667,           0.0416667,           0.0833334,  502.1450692400854, 0.3,  0.4,  0.11, 0.37),
    @(15, 16, 286.7480913238156,  0.0138889,           0.0138889,           0.0277778,  430.1221369857234, 0.33, 0.49, 0.12, 0.36),
    @(16, 17, 262.4045378156197,  0,                   0,                   0,          393.6068067234295, 0.39, 0.46, 0.14, 0.39),
    @(17, 18, 309.1354664251031,  0,                   0,                   0,          463.7031996376547, 0.37, 0.48, 0.19, 0.34),
    @(18, 19, 384.636049180636,   0,                   0,                   0,          576.954073770954,  0.35, 0.4,  0.18, 0.39),
    @(19, 20, 433.5207499755035,  0,                   0,                   0,          650.2811249632553, 0.38, 0.41, 0.17, 0.33),
    @(20, 21, 393.8044005019047,  0,                   0,                   0,          590.7066007528571, 0.35, 0.49, 0.16, 0.37),
    @(21, 22, 333.3011855326704,  0,                   0,                   0,          499.9517782990056, 0.3,  0.4,  0.19, 0.34),
    @(22, 23, 282.2034344188766,  0,                   0,                   0,          423.3051516283149, 0.39, 0.41, 0.14, 0.3),
    @(23, 24, 208.1255268683671,  0,                   0,                   0,          312.1882903025506, 0.32, 0.45, 0.18, 0.39),
    @(24, 25, 165.1488800499209,  0,                   0,                   0,          247.7233200748814, 0.4,  0.48, 0.2,  0.4)
)

$startRow = 2
$lastOldRow = 12
$lastNewRow = $startRow + $data.Count - 1   # 17

# Write the new values across the full block A2:K17
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# New rows below the old table end (13..17) need the same column-A formatting
# (bold/bordered/centered style) as the existing data rows - copy it down.
if ($lastNewRow -gt $lastOldRow) {
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range("A" + ($lastOldRow + 1) + ":A" + $lastNewRow).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
